$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.47"
$ws.Range("E2").Value = "'-4.74%"
$ws.Range("D3").Value = "'39.89"
$ws.Range("E3").Value = "'-7.43%"
$ws.Range("D4").Value = "'5.121"
$ws.Range("E4").Value = "'-1.53%"
$ws.Range("D5").Value = "'0.07708"
$ws.Range("E5").Value = "'-5.85%"
$ws.Range("D6").Value = "'4.238"
$ws.Range("E6").Value = "'-1.91%"
$ws.Range("D7").Value = "'1.625"
$ws.Range("E7").Value = "'-11.45%"
$ws.Range("D8").Value = "'0.8804"
$ws.Range("E8").Value = "'-5.96%"
$ws.Range("D9").Value = "'0.1004"
$ws.Range("E9").Value = "'-9.92%"
$ws.Range("D10").Value = "'0.1747"
$ws.Range("E10").Value = "'-6.27%"
$ws.Range("D11").Value = "'0.08911"
$ws.Range("E11").Value = "'-4.60%"
$ws.Range("D12").Value = "'0.04413"
$ws.Range("E12").Value = "'-4.53%"
$ws.Range("D13").Value = "'0.1057"
$ws.Range("E13").Value = "'0.01%"
$ws.Range("D14").Value = "'0.001259"
$ws.Range("E14").Value = "'-2.68%"
$ws.Range("D15").Value = "'0.005894"
$ws.Range("E15").Value = "'1.91%"
$ws.Range("D16").Value = "'3.353"
$ws.Range("E16").Value = "'-0.10%"
$ws.Range("D17").Value = "'2.436"
$ws.Range("E17").Value = "'-2.99%"
$ws.Range("D18").Value = "'0.3323"
$ws.Range("E18").Value = "'-0.48%"
$ws.Range("D19").Value = "'7.034"
$ws.Range("E19").Value = "'-5.06%"
$ws.Range("E20").Value = "'-4.06%"
$ws.Range("E21").Value = "'14.42%"
$ws.Range("D22").Value = "'0.04144"
$ws.Range("E22").Value = "'-0.15%"
$ws.Range("D23").Value = "'0.001201"
$ws.Range("E23").Value = "'-3.88%"
$ws.Range("D24").Value = "'0.004082"
$ws.Range("E24").Value = "'-5.28%"
$ws.Range("E25").Value = "'10.91%"
$ws.Range("E26").Value = "'0.17%"
$ws.Range("D38").Value = "'0.02337"
$ws.Range("E38").Value = "'-14.30%"
$ws.Range("D39").Value = "'0.05144"
$ws.Range("E39").Value = "'-7.17%"
$ws.Range("D40").Value = "'0.007926"
$ws.Range("E40").Value = "'-0.61%"
$ws.Range("D41").Value = "'0.1324"
$ws.Range("E41").Value = "'-4.99%"
$ws.Range("D42").Value = "'0.006336"
$ws.Range("E42").Value = "'-3.00%"
$ws.Range("D43").Value = "'0.001943"
$ws.Range("E43").Value = "'-7.13%"
$ws.Range("D44").Value = "'0.008585"
$ws.Range("E44").Value = "'15.10%"
$ws.Range("D45").Value = "'0.3053"
$ws.Range("E45").Value = "'-4.75%"
$ws.Range("D46").Value = "'0.00006514"
$ws.Range("E46").Value = "'-6.53%"
$ws.Range("E47").Value = "'0.10%"
$ws.Range("D48").Value = "'0.007002"
$ws.Range("E48").Value = "'98.53%"
$ws.Range("D49").Value = "'0.002184"
$ws.Range("E49").Value = "'-36.90%"
$ws.Range("E50").Value = "'0.10%"
$ws.Range("E51").Value = "'0.10%"
